# Generate Report for Archive
# - Update status text "Ready for handoff" -> "In Translation" on all sheets
# - Shrink the now-narrower "Status"/zh-cn/de-de columns to fit the new text

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status text ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: column C holds the Status text ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: column C holds the Status text ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
